$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 19125
$ws.Range("B3").Value = 14331
$ws.Range("B4").Value = 1802
$ws.Range("B5").Value = 17671
